$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 359.5
$ws.Range("I2").Value = 409.875
$ws.Range("J2").Value = 258.75
$ws.Range("K2").Value = 409.875
$ws.Range("L2").Value = 258.75
$ws.Range("M2").Value = -296.875
$ws.Range("N2").Value = -484.75
$ws.Range("H69").Value = 9659.615
$ws.Range("I69").Value = 7132.8887
$ws.Range("J69").Value = 15344.75
$ws.Range("K69").Value = 21398.6661
$ws.Range("L69").Value = 46034.25
$ws.Range("M69").Value = -20524.6661
$ws.Range("N69").Value = -47782.25
$ws.Range("H72").Value = 9659.615
$ws.Range("I72").Value = 7132.8887
$ws.Range("J72").Value = 15344.75
$ws.Range("K72").Value = 64195.99830000001
$ws.Range("L72").Value = 138102.75
$ws.Range("M72").Value = -59827.99830000001
$ws.Range("N72").Value = -146838.75
$ws.Range("H75").Value = 152578.5
$ws.Range("J75").Value = 196771.33
$ws.Range("L75").Value = 196771.33
$ws.Range("N75").Value = -198643.33
$ws.Range("H78").Value = 152578.5
$ws.Range("J78").Value = 196771.33
$ws.Range("L78").Value = 590313.99
$ws.Range("N78").Value = -599673.99
$ws.Range("H80").Value = 56438.785
$ws.Range("I80").Value = 92132
$ws.Range("J80").Value = 4088.7334
$ws.Range("K80").Value = 276396
$ws.Range("L80").Value = 12266.2002
$ws.Range("M80").Value = -275398
$ws.Range("N80").Value = -14262.2002
$ws.Range("H83").Value = 56438.785
$ws.Range("I83").Value = 92132
$ws.Range("J83").Value = 4088.7334
$ws.Range("K83").Value = 829188
$ws.Range("L83").Value = 36798.6006
$ws.Range("M83").Value = -824196
$ws.Range("N83").Value = -46782.6006
$ws.Range("H86").Value = 26137.482
$ws.Range("I86").Value = 5011.5
$ws.Range("K86").Value = 5011.5
$ws.Range("M86").Value = -3888.5
$ws.Range("H89").Value = 26137.482
$ws.Range("I89").Value = 5011.5
$ws.Range("K89").Value = 25057.5
$ws.Range("M89").Value = -19441.5
$ws.Range("H132").Value = 4353772.5
$ws.Range("I132").Value = 6238.2
$ws.Range("K132").Value = 18714.6
$ws.Range("M132").Value = -16184.6
$ws.Range("H141").Value = 7264.033
$ws.Range("I141").Value = 7212.385
$ws.Range("K141").Value = 21637.155
$ws.Range("M141").Value = -16457.155

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 49994.332
$ws.Range("J141").Value = 49994.332
$ws.Range("L141").Value = 49994.332
$ws.Range("N141").Value = -60354.332

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4370.1665
$ws.Range("I86").Value = 4963.5386
$ws.Range("J86").Value = 2827.4
$ws.Range("K86").Value = 4963.5386
$ws.Range("L86").Value = 2827.4
$ws.Range("M86").Value = -3840.5386
$ws.Range("N86").Value = -5073.4
$ws.Range("H89").Value = 4370.1665
$ws.Range("I89").Value = 4963.5386
$ws.Range("J89").Value = 2827.4
$ws.Range("K89").Value = 24817.693
$ws.Range("L89").Value = 14137
$ws.Range("M89").Value = -19201.693
$ws.Range("N89").Value = -25369
$ws.Range("H94").Value = 4005.5386
$ws.Range("I94").Value = 2791.2
$ws.Range("K94").Value = 2791.2
$ws.Range("M94").Value = -2340.2
$ws.Range("H132").Value = 72779.8
$ws.Range("J132").Value = 72779.8
$ws.Range("L132").Value = 72779.8
$ws.Range("N132").Value = -82899.8
$ws.Range("H134").Value = 11805.8
$ws.Range("I134").Value = 21681.5
$ws.Range("K134").Value = 65044.5
$ws.Range("M134").Value = -62509.5
$ws.Range("H138").Value = 86099.89999999999
$ws.Range("J138").Value = 90111
$ws.Range("L138").Value = 90111
$ws.Range("N138").Value = -100391

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 85165.25
$ws.Range("I16").Value = 1756
$ws.Range("J16").Value = 168574.5
$ws.Range("K16").Value = 1756
$ws.Range("L16").Value = 168574.5
$ws.Range("M16").Value = -1469
$ws.Range("N16").Value = -169148.5
$ws.Range("H31").Value = 10051.952
$ws.Range("J31").Value = 4370
$ws.Range("L31").Value = 4370
$ws.Range("N31").Value = -4960
$ws.Range("H34").Value = 10051.952
$ws.Range("J34").Value = 4370
$ws.Range("L34").Value = 4370
$ws.Range("N34").Value = -4774
$ws.Range("H105").Value = 6395.9473
$ws.Range("I105").Value = 8196.071
$ws.Range("K105").Value = 8196.071
$ws.Range("M105").Value = -6449.071
$ws.Range("H113").Value = 85165.25
$ws.Range("I113").Value = 1756
$ws.Range("J113").Value = 168574.5
$ws.Range("K113").Value = 1756
$ws.Range("L113").Value = 168574.5
$ws.Range("M113").Value = 414
$ws.Range("N113").Value = -172914.5
$ws.Range("H132").Value = 1358.2333
$ws.Range("I132").Value = 1311.8928
$ws.Range("K132").Value = 3935.6784
$ws.Range("M132").Value = -1405.6784
$ws.Range("H134").Value = 6489
$ws.Range("I134").Value = 5491.3335
$ws.Range("K134").Value = 16474.0005
$ws.Range("M134").Value = -13939.0005
$ws.Range("H141").Value = 102141.12
$ws.Range("J141").Value = 105029.305
$ws.Range("L141").Value = 105029.305
$ws.Range("N141").Value = -115389.305

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 98279
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 98279
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 294837
$ws.Range("N80").Value = -296709
$ws.Range("H83").Value = 98279
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 98279
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 884511
$ws.Range("N83").Value = -893871
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14645
$ws.Range("I80").Value = 14645
$ws.Range("K80").Value = 14645
$ws.Range("M80").Value = -13647
$ws.Range("H83").Value = 14645
$ws.Range("I83").Value = 14645
$ws.Range("K83").Value = 73225
$ws.Range("M83").Value = -68233
$ws.Range("H113").Value = 7487.5264
$ws.Range("I113").Value = 8203.9375
$ws.Range("K113").Value = 8203.9375
$ws.Range("M113").Value = -6033.9375
$ws.Range("H122").Value = 7301.8623
$ws.Range("I122").Value = 4768.6665
$ws.Range("J122").Value = 41500
$ws.Range("K122").Value = 14305.9995
$ws.Range("L122").Value = 124500
$ws.Range("M122").Value = -11855.9995
$ws.Range("N122").Value = -129400

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1187.1489
$ws.Range("I16").Value = 1093.1945
$ws.Range("K16").Value = 1093.1945
$ws.Range("M16").Value = -923.1945000000001
$ws.Range("H55").Value = 1207.6111
$ws.Range("I55").Value = 430.63635
$ws.Range("K55").Value = 430.63635
$ws.Range("M55").Value = -257.63635
$ws.Range("H61").Value = 5854.125
$ws.Range("I61").Value = 2555.5833
$ws.Range("K61").Value = 2555.5833
$ws.Range("M61").Value = -2353.5833
$ws.Range("H74").Value = 27999.2
$ws.Range("I74").Value = 27999.2
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 27999.2
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -27001.2
$ws.Range("H77").Value = 27999.2
$ws.Range("I77").Value = 27999.2
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 83997.60000000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -79005.60000000001
$ws.Range("H113").Value = 5854.125
$ws.Range("I113").Value = 2555.5833
$ws.Range("K113").Value = 2555.5833
$ws.Range("M113").Value = -385.5832999999998
$ws.Range("H122").Value = 3934.2144
$ws.Range("I122").Value = 3558.6487
$ws.Range("K122").Value = 10675.9461
$ws.Range("M122").Value = -8225.946100000001
$ws.Range("H132").Value = 598310.3
$ws.Range("I132").Value = 877965.9
$ws.Range("J132").Value = 4042.25
$ws.Range("K132").Value = 2633897.7
$ws.Range("L132").Value = 12126.75
$ws.Range("M132").Value = -2631367.7
$ws.Range("N132").Value = -17186.75
$ws.Range("H136").Value = 4960.1284
$ws.Range("I136").Value = 3750
$ws.Range("J136").Value = 5716.4585
$ws.Range("K136").Value = 11250
$ws.Range("L136").Value = 17149.3755
$ws.Range("M136").Value = -8700
$ws.Range("N136").Value = -22249.3755
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 5000
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 5000
$ws.Range("N25").Value = -5586
$ws.Range("H122").Value = 10746.017
$ws.Range("I122").Value = 1713.7234
$ws.Range("K122").Value = 5141.1702
$ws.Range("M122").Value = -2691.1702
$ws.Range("H132").Value = 12255.718
$ws.Range("I132").Value = 14145.5
$ws.Range("K132").Value = 42436.5
$ws.Range("M132").Value = -39906.5

Write-Host "Applied all cell updates"